# Updates the cryptocurrency Price (D) and Volume(1h) (E) columns on Sheet1
# to the latest scraped values, matching the GitHub Actions refresh job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.078.54'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').Value = '1.621.48'
$ws.Range('E3').Value = '  -0.97%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = '''214.07'
$ws.Range('E5').Value = '  -1.18%  '
$ws.Range('D6').Value = '''0.515'
$ws.Range('E6').Value = '  -0.29%  '
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('D8').Value = '''0.0628'
$ws.Range('E8').Value = '  +0.56%  '
$ws.Range('E9').Value = '  -1.31%  '
$ws.Range('D10').Value = '''19.90'
$ws.Range('E10').Value = '  +0.05%  '
$ws.Range('E11').Value = '  -0.48%  '
$ws.Range('E12').Value = '  -0.99%  '
$ws.Range('D13').Value = '1.621.53'
$ws.Range('E13').Value = '  -1.00%  '
$ws.Range('D14').Value = '''4.12'
$ws.Range('E14').Value = '  +0.11%  '
$ws.Range('D15').Value = '''0.538'
$ws.Range('E15').Value = '  -0.34%  '
$ws.Range('D16').Value = '27.038.23'
$ws.Range('E16').Value = '  -0.24%  '
$ws.Range('D17').Value = '''64.48'
$ws.Range('E17').Value = '  -3.17%  '
$ws.Range('E18').Value = '  -0.15%  '
$ws.Range('D19').Value = '''214.23'
$ws.Range('E19').Value = '  -1.16%  '
$ws.Range('E20').Value = '  -0.04%  '
$ws.Range('E21').Value = '  -0.59%  '
$ws.Range('D22').Value = '''4.34'
$ws.Range('E22').Value = '  -1.60%  '
$ws.Range('E23').Value = '  -7.56%  '
$ws.Range('D24').Value = '''9.01'
$ws.Range('E24').Value = '  -0.88%  '
$ws.Range('D25').Value = '''147.80'
$ws.Range('E25').Value = '  +0.67%  '
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('D27').Value = '''7.40'
$ws.Range('E27').Value = '  +0.11%  '
$ws.Range('E28').Value = '  -3.26%  '
$ws.Range('D29').Value = '''15.50'
$ws.Range('E29').Value = '  -1.02%  '
$ws.Range('E30').Value = '  +0.51%  '
$ws.Range('E31').Value = '  -0.99%  '
$ws.Range('D32').Value = '''3.32'
$ws.Range('E32').Value = '  -1.73%  '
$ws.Range('D33').Value = '''0.704'
$ws.Range('E33').Value = '  +30.07%  '
$ws.Range('D34').Value = '''2.99'
$ws.Range('E34').Value = '  -0.34%  '
$ws.Range('D35').Value = '1.342.93'
$ws.Range('E35').Value = '  +3.19%  '
$ws.Range('E36').Value = '  -0.74%  '
$ws.Range('E37').Value = '  -0.71%  '
$ws.Range('E38').Value = '  -0.03%  '
$ws.Range('D39').Value = '''0.841'
$ws.Range('E39').Value = '  -1.50%  '
$ws.Range('E40').Value = '  -0.14%  '
$ws.Range('E41').Value = '  +0.45%  '
$ws.Range('D42').Value = '''0.794'
$ws.Range('E42').Value = '  -1.69%  '
$ws.Range('D43').Value = '''5.32'
$ws.Range('E43').Value = '  +0.42%  '
$ws.Range('D44').Value = '''63.75'
$ws.Range('E44').Value = '  +3.48%  '
$ws.Range('D45').Value = '1.759.49'
$ws.Range('E45').Value = '  -1.02%  '
$ws.Range('D46').Value = '''89.90'
$ws.Range('E46').Value = '  -1.39%  '
$ws.Range('D47').Value = '''1.65'
$ws.Range('E47').Value = '  +3.00%  '
$ws.Range('D48').Value = '''0.852'
$ws.Range('E48').Value = '  +27.25%  '
$ws.Range('D49').Value = '''0.100'
$ws.Range('E49').Value = '  +4.85%  '
$ws.Range('E50').Value = '  -0.02%  '
$ws.Range('D51').Value = '''7.57'
$ws.Range('E51').Value = '  -0.80%  '
